$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '43.401.68'
$ws.Range("E2").Value = '  -0.52%  '
$ws.Range("D3").Value = '2.281.62'
$ws.Range("E3").Value = '  -0.52%  '
$ws.Range("E4").Value = '  -0.04%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '113.32'
$ws.Range("E5").Value = '  -1.34%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '265.16'
$ws.Range("E6").Value = '  -1.59%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.622'
$ws.Range("E7").Value = '  -0.32%  '
$ws.Range("E8").Value = '  +0.06%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.606'
$ws.Range("E9").Value = '  -1.58%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '47.56'
$ws.Range("E10").Value = '  -0.90%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0934'
$ws.Range("E11").Value = '  -0.65%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '9.29'
$ws.Range("E12").Value = '  +8.12%  '
$ws.Range("E13").Value = '  +0.98%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '15.47'
$ws.Range("E14").Value = '  -0.91%  '
$ws.Range("D15").Value = '2.608.06'
$ws.Range("E15").Value = '  -1.08%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.867'
$ws.Range("E16").Value = '  +1.98%  '
$ws.Range("D17").Value = '2.273.48'
$ws.Range("E17").Value = '  -0.82%  '
$ws.Range("D18").Value = '43.208.33'
$ws.Range("E18").Value = '  -0.98%  '
$ws.Range("E19").Value = '  -1.23%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '6.81'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '71.91'
$ws.Range("E21").Value = '  -0.93%  '
$ws.Range("E22").Value = '  -1.28%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '233.95'
$ws.Range("E23").Value = '  +0.16%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '9.61'
$ws.Range("E24").Value = '  +0.42%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.88'
$ws.Range("E25").Value = '  +0.98%  '
$ws.Range("E26").Value = '  +1.66%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '11.42'
$ws.Range("E27").Value = '  -1.21%  '
$ws.Range("E28").Value = '  +0.20%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '40.91'
$ws.Range("E29").Value = '  -3.00%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '3.34'
$ws.Range("E30").Value = '  -2.05%  '
$ws.Range("E31").Value = '  -0.83%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '173.75'
$ws.Range("E32").Value = '  -1.66%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '21.53'
$ws.Range("E33").Value = '  -0.55%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.0906'
$ws.Range("E34").Value = '  -2.03%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '5.71'
$ws.Range("E35").Value = '  +3.03%  '
$ws.Range("E36").Value = '  +0.50%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '4.66'
$ws.Range("E37").Value = '  -1.40%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.0366'
$ws.Range("E38").Value = '  +2.37%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.89'
$ws.Range("E39").Value = '  +1.44%  '
$ws.Range("E40").Value = '  -5.05%  '
$ws.Range("E41").Value = '  +8.88%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '76.48'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '14.17'
$ws.Range("E43").Value = '  +2.44%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.239'
$ws.Range("E44").Value = '  -1.69%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '6.22'
$ws.Range("E45").Value = '  +4.01%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.999'
$ws.Range("E46").Value = '  -0.08%  '
$ws.Range("E47").Value = '  -3.74%  '
$ws.Range("E48").Value = '  -1.47%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '103.87'
$ws.Range("E49").Value = '  +0.98%  '
$ws.Range("E50").Value = '  +0.75%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0996'
$ws.Range("E51").Value = '  -0.67%  '
